$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(813,     803.85, 808.65,  810.5,   42,  807)
    3  = @(872.4,   853.5,  869,     868.65,  54,  856)
    4  = @(44759.95,44470.05,44689.3,44644.05,30,  44663)
    5  = @(340.25,  334.1,  336.15,  335.9,   122, 338.75)
    6  = @(524.9,   515.7,  521,     521.4,   70,  522.75)
    7  = @(482.05,  476,    479.65,  479.8,   75,  481.1)
    8  = @(972.05,  956.85, 964.5,   962.05,  254, 970.2)
    9  = @(707.45,  699,    706.1,   706.6,   52,  703.95)
    10 = @(19673,   19556.1,19665.2, 19658.85,54,  19620.05)
    11 = @(2447.9,  2431,   2438.25, 2439.45, 53,  2440.45)
    12 = @(577,     571.5,  574.2,   573.65,  219, 575.25)
    13 = @(887.25,  850.6,  886.7,   882.55,  104, 851.35)
    14 = @(614.65,  607.9,  613.65,  613.5,   102, 611.7)
    15 = @(131.5,   128.65, 129.95,  130,     455, 131.3)
    16 = @(3454.65, 3414.4, 3437,    3440.65, 15,  3453.1)
    17 = @(3185.35, 3141,   3181.35, 3181.25, 14,  3147.3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
